$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.899.44'
$ws.Range('E2').Value = '  -3.94%  '
$ws.Range('D3').Value = '3.304.67'
$ws.Range('E3').Value = '  -4.24%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.26%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '557.09'
$c.ClearFormats()
$ws.Range('E5').Value = '  -2.49%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '185.59'
$c.ClearFormats()
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  -4.30%  '
$ws.Range('D9').Value = '3.295.68'
$ws.Range('E9').Value = '  -4.14%  '
$ws.Range('E10').Value = '  -7.65%  '
$ws.Range('E11').Value = '  -4.51%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '47.62'
$c.ClearFormats()
$ws.Range('E12').Value = '  -3.44%  '
$ws.Range('E13').Value = '  -3.64%  '
$ws.Range('E14').Value = '  -3.54%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '633.57'
$c.ClearFormats()
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Value = '3.834.18'
$ws.Range('E16').Value = '  -4.74%  '
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '65.914.48'
$ws.Range('E18').Value = '  -3.76%  '
$ws.Range('E19').Value = '  -2.59%  '
$ws.Range('D20').Value = '3.297.91'
$ws.Range('E20').Value = '  -5.28%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '11.38'
$c.ClearFormats()
$ws.Range('E21').Value = '  -5.92%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.907'
$c.ClearFormats()
$ws.Range('E22').Value = '  -2.91%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '18.06'
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.15%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '103.17'
$c.ClearFormats()
$ws.Range('E24').Value = '  +6.17%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.96'
$c.ClearFormats()
$ws.Range('E25').Value = '  -4.28%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '3.95'
$c.ClearFormats()
$ws.Range('E26').Value = '  -6.21%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '5.94'
$c.ClearFormats()
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('E28').Value = '  -4.16%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.59'
$c.ClearFormats()
$ws.Range('E29').Value = '  -1.62%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '8.67'
$c.ClearFormats()
$ws.Range('E30').Value = '  -4.89%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '30.21'
$c.ClearFormats()
$ws.Range('E31').Value = '  -4.98%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.05'
$c.ClearFormats()
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('E33').Value = '  -2.23%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '11.12'
$c.ClearFormats()
$ws.Range('E34').Value = '  -2.55%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '542.36'
$c.ClearFormats()
$ws.Range('E35').Value = '  -4.86%  '
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('D37').Value = '3.808.83'
$ws.Range('E37').Value = '  -1.47%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '57.56'
$c.ClearFormats()
$ws.Range('E38').Value = '  -2.98%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '0.0₃0738'
$ws.Range('E40').Value = '  -3.91%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '33.94'
$c.ClearFormats()
$ws.Range('E41').Value = '  +2.63%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.129'
$c.ClearFormats()
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.70'
$c.ClearFormats()
$ws.Range('E43').Value = '  -4.04%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '3.25'
$c.ClearFormats()
$ws.Range('E44').Value = '  -6.09%  '
$ws.Range('B45').Value = 'CoreDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.24'
$c.ClearFormats()
$ws.Range('E45').Value = '  -13.85%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.335'
$c.ClearFormats()
$ws.Range('E46').Value = '  -7.79%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0418'
$c.ClearFormats()
$ws.Range('E47').Value = '  -3.30%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.26'
$c.ClearFormats()
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.60'
$c.ClearFormats()
$ws.Range('E49').Value = '  -6.49%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.129'
$c.ClearFormats()
$ws.Range('E50').Value = '  -3.99%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E51').Value = '  -0.20%  '
